# Update the NATMI TPM output values (new TPM recompute) in the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("M2").Value = 0.162165
$ws.Range("N2").Value = 0.486495
$ws.Range("O2").Value = 0.1006291402646046
$ws.Range("P2").Value = 0.1006291402646046
$ws.Range("Q2").Value = 1.48493831029
$ws.Range("R2").Value = 13.36444479261
$ws.Range("S2").Value = 0.09754703994695904
$ws.Range("T2").Value = 0.09754703994695905

# --- Row 3 ---
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.1908661724170313
$ws.Range("P3").Value = 0.1908661724170313
$ws.Range("S3").Value = 0.1850202644714057
$ws.Range("T3").Value = 0.1850202644714057

# --- Row 4 ---
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 1.141763333333333
$ws.Range("N4").Value = 3.42529
$ws.Range("O4").Value = 0.7085046873183641
$ws.Range("P4").Value = 0.7085046873183641
$ws.Range("Q4").Value = 10.45508041162445
$ws.Range("R4").Value = 94.09572370462001
$ws.Range("S4").Value = 0.6868043874241656
$ws.Range("T4").Value = 0.6868043874241656

# --- Row 5 ---
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("M5").Value = 0.162165
$ws.Range("N5").Value = 0.486495
$ws.Range("O5").Value = 0.1006291402646046
$ws.Range("P5").Value = 0.1006291402646046
$ws.Range("Q5").Value = 0.046918172405
$ws.Range("R5").Value = 0.422263551645
$ws.Range("S5").Value = 0.003082100317645544
$ws.Range("T5").Value = 0.003082100317645544

# --- Row 6 ---
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.1908661724170313
$ws.Range("P6").Value = 0.1908661724170313
$ws.Range("Q6").Value = 0.08899104136433335
$ws.Range("R6").Value = 0.8009193722790001
$ws.Range("S6").Value = 0.005845907945625563
$ws.Range("T6").Value = 0.005845907945625563

# --- Row 7 ---
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 1.141763333333333
$ws.Range("N7").Value = 3.42529
$ws.Range("O7").Value = 0.7085046873183641
$ws.Range("P7").Value = 0.7085046873183641
$ws.Range("Q7").Value = 0.3303391540655556
$ws.Range("R7").Value = 2.97305238659
$ws.Range("S7").Value = 0.02170029989419852
$ws.Range("T7").Value = 0.02170029989419852
